$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "未平倉績效" column (old B),
# shifting old B,C,D -> D,E,F
$ws.Range("B:C").Insert()

# Insert one more new column before the old "已平倉績效" column
# (which is now at E after the first insert), shifting it and the
# last column right by one: E,F -> F,G
$ws.Range("E:E").Insert()

# Fill in the new header cells (row 1)
$ws.Range("B1").Value = "未平倉成本"
$ws.Range("C1").Value = "未實現損益"
$ws.Range("E1").Value = "已實現損益"

# The cells left behind in the two freshly inserted columns (B,C,E) for
# the data rows should remain blank, just like their neighbouring
# (pre-existing) blank cells. Copy one of those already-blank cells
# over them so they are materialized the same way instead of being
# left as completely empty/absent cells.
$ws.Range("F2").Copy($ws.Range("B2:C5"))
$ws.Range("F2").Copy($ws.Range("E2:E5"))

# Fill in the new values for row 6 (2021-12-30), preserving them as
# literal text (not auto-converted to numbers) by forcing a Text
# number format before assignment.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "56348.0"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "-318.0"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-363.0"
